$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "Drinks" block of rows 61-68 was manually reordered in Excel (rows were
# dragged/moved around - not a plain ascending sort, since the final two
# rows end up as 35 then 30). Re-create the exact resulting row contents.
# Columns A, C and F are identical for every row in this block, so only
# B (item), D (price) and E (image) need to be rewritten.

$ws.Range("B61").Value = "AppyFizz20"
$ws.Range("D61").Value = 20
$ws.Range("E61").Value = "AppyFizz20.jpg"

$ws.Range("B62").Value = "Cream bell kesar badam milk"
$ws.Range("D62").Value = 20
$ws.Range("E62").Value = "Cream bell kesar badam milk.jpg"

$ws.Range("B63").Value = "DEW20"
$ws.Range("D63").Value = 20
$ws.Range("E63").Value = "Dew_20.jpg"

$ws.Range("B64").Value = "Sting 20 rs"
$ws.Range("D64").Value = 20
$ws.Range("E64").Value = "Sting 20 rs.jpg"

$ws.Range("B65").Value = "Dite Coke 25 rs"
$ws.Range("D65").Value = 25
$ws.Range("E65").Value = "Dite Coke 25 rs.jpg"

$ws.Range("B66").Value = "Thumsup30"
$ws.Range("D66").Value = 30
$ws.Range("E66").Value = "Thusmup Can.jpg"

$ws.Range("B67").Value = "Thums up 35"
$ws.Range("D67").Value = 35
$ws.Range("E67").Value = "Thums up 35.jpg"

$ws.Range("B68").Value = "Campa Energy 30Rs"
$ws.Range("D68").Value = 30
$ws.Range("E68").Value = "Campa Energy 30Rs.jpg"

# The AutoFilter-derived named range shrank along with the reorder/cleanup.
$wb.Names.Item("Sheet1!_FilterDatabase").RefersTo = "=Sheet1!`$A`$1:`$F`$61"

# Leave the selection where the author ended up.
$ws.Range("B68").Select()
